$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Row 7: comment for Lab 7 (Branch&Bound) column (added to shared strings first)
$ws.Range("H7").Value = "Not done"

# Row 6: Lab 7 (Branch&Bound) mark = 0, Test mark = "Passed"
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = "Passed"

# Update view: scroll so column D is the left-most visible column, and select I6
[void]$ws.Range("I6").Select()
$excel.ActiveWindow.ScrollColumn = 4
